# Update "想去人数" (interest count) figures for several events.
# Sheet "展览" (sheet1) and sheet "全部类型" (sheet4) both list the same
# events; each gets the corresponding F-column value bumped.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 225
$ws1.Range("F4").Value = 13010
$ws1.Range("F6").Value = 217
$ws1.Range("F10").Value = 227
$ws1.Range("F18").Value = 5540
$ws1.Range("F19").Value = 107
$ws1.Range("F20").Value = 55
$ws1.Range("F24").Value = 141

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 225
$ws4.Range("F4").Value = 13011
$ws4.Range("F6").Value = 217
$ws4.Range("F10").Value = 227
$ws4.Range("F18").Value = 5540
$ws4.Range("F19").Value = 107
$ws4.Range("F20").Value = 55
$ws4.Range("F24").Value = 141
